$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Rename sheets: Sheet1 -> A01, Sheet2 -> GA02 ---
$ws1.Name = "A01"
$ws2.Name = "GA02"

# --- Values & formulas for GA02 (Sheet2) ---
$ws2.Range("A1").Value = "Week 1"
$ws2.Range("B1").Value = "GA02"
$ws2.Range("I1").Value = "Week 2"
$ws2.Range("J1").Value = "GA02"
$ws2.Range("C3").Value = "Documentation"
$ws2.Range("K3").Value = "Documentation"
$ws2.Range("B4").Value = "Måndag"
$ws2.Range("C4").Value = "Tisdag"
$ws2.Range("D4").Value = "Onsdag"
$ws2.Range("E4").Value = "Torsdag"
$ws2.Range("F4").Value = "Fredag"
$ws2.Range("J4").Value = "Måndag"
$ws2.Range("K4").Value = "Tisdag"
$ws2.Range("L4").Value = "Onsdag"
$ws2.Range("M4").Value = "Torsdag"
$ws2.Range("N4").Value = "Fredag"
$ws2.Range("A5").Value = "Calle"
$ws2.Range("G5").Formula = "=SUM(B5:F5)"
$ws2.Range("I5").Value = "Calle"
$ws2.Range("O5").Formula = "=SUM(J5:N5)"
$ws2.Range("A6").Value = "Kim"
$ws2.Range("G6").Formula = "=SUM(B6:F6)"
$ws2.Range("I6").Value = "Kim"
$ws2.Range("O6").Formula = "=SUM(J6:N6)"
$ws2.Range("A7").Value = "Nils"
$ws2.Range("G7").Formula = "=SUM(B7:F7)"
$ws2.Range("I7").Value = "Nils"
$ws2.Range("O7").Formula = "=SUM(J7:N7)"
$ws2.Range("A8").Value = "Rasmus"
$ws2.Range("G8").Formula = "=SUM(B8:F8)"
$ws2.Range("I8").Value = "Rasmus"
$ws2.Range("O8").Formula = "=SUM(J8:N8)"
$ws2.Range("G9").Formula = "=SUM(G5:G8)"
$ws2.Range("O9").Formula = "=SUM(O5:O8)"
$ws2.Range("C10").Value = "Strategy meeting"
$ws2.Range("K10").Value = "Strategy meeting"
$ws2.Range("B11").Value = "Måndag"
$ws2.Range("C11").Value = "Tisdag"
$ws2.Range("D11").Value = "Onsdag"
$ws2.Range("E11").Value = "Torsdag"
$ws2.Range("F11").Value = "Fredag"
$ws2.Range("J11").Value = "Måndag"
$ws2.Range("K11").Value = "Tisdag"
$ws2.Range("L11").Value = "Onsdag"
$ws2.Range("M11").Value = "Torsdag"
$ws2.Range("N11").Value = "Fredag"
$ws2.Range("A12").Value = "Calle"
$ws2.Range("G12").Formula = "=SUM(B12:F12)"
$ws2.Range("I12").Value = "Calle"
$ws2.Range("O12").Formula = "=SUM(J12:N12)"
$ws2.Range("A13").Value = "Kim"
$ws2.Range("G13").Formula = "=SUM(B13:F13)"
$ws2.Range("I13").Value = "Kim"
$ws2.Range("O13").Formula = "=SUM(J13:N13)"
$ws2.Range("A14").Value = "Nils"
$ws2.Range("G14").Formula = "=SUM(B14:F14)"
$ws2.Range("I14").Value = "Nils"
$ws2.Range("O14").Formula = "=SUM(J14:N14)"
$ws2.Range("A15").Value = "Rasmus"
$ws2.Range("G15").Formula = "=SUM(B15:F15)"
$ws2.Range("I15").Value = "Rasmus"
$ws2.Range("O15").Formula = "=SUM(J15:N15)"
$ws2.Range("G16").Formula = "=SUM(G12:G15)"
$ws2.Range("O16").Formula = "=SUM(O12:O15)"
$ws2.Range("C17").Value = "Updating A01"
$ws2.Range("K17").Value = "Updating A01"
$ws2.Range("B18").Value = "Måndag"
$ws2.Range("C18").Value = "Tisdag"
$ws2.Range("D18").Value = "Onsdag"
$ws2.Range("E18").Value = "Torsdag"
$ws2.Range("F18").Value = "Fredag"
$ws2.Range("J18").Value = "Måndag"
$ws2.Range("K18").Value = "Tisdag"
$ws2.Range("L18").Value = "Onsdag"
$ws2.Range("M18").Value = "Torsdag"
$ws2.Range("N18").Value = "Fredag"
$ws2.Range("A19").Value = "Calle"
$ws2.Range("G19").Formula = "=SUM(B19:F19)"
$ws2.Range("I19").Value = "Calle"
$ws2.Range("K19").Value = 4
$ws2.Range("O19").Formula = "=SUM(J19:N19)"
$ws2.Range("A20").Value = "Kim"
$ws2.Range("G20").Formula = "=SUM(B20:F20)"
$ws2.Range("I20").Value = "Kim"
$ws2.Range("O20").Formula = "=SUM(J20:N20)"
$ws2.Range("A21").Value = "Nils"
$ws2.Range("G21").Formula = "=SUM(B21:F21)"
$ws2.Range("I21").Value = "Nils"
$ws2.Range("O21").Formula = "=SUM(J21:N21)"
$ws2.Range("A22").Value = "Rasmus"
$ws2.Range("G22").Formula = "=SUM(B22:F22)"
$ws2.Range("I22").Value = "Rasmus"
$ws2.Range("K22").Value = 4
$ws2.Range("O22").Formula = "=SUM(J22:N22)"
$ws2.Range("G23").Formula = "=SUM(G19:G22)"
$ws2.Range("O23").Formula = "=SUM(O19:O22)"
$ws2.Range("C24").Value = "Assignment planning"
$ws2.Range("K24").Value = "Assignment planning"
$ws2.Range("B25").Value = "Måndag"
$ws2.Range("C25").Value = "Tisdag"
$ws2.Range("D25").Value = "Onsdag"
$ws2.Range("E25").Value = "Torsdag"
$ws2.Range("F25").Value = "Fredag"
$ws2.Range("J25").Value = "Måndag"
$ws2.Range("K25").Value = "Tisdag"
$ws2.Range("L25").Value = "Onsdag"
$ws2.Range("M25").Value = "Torsdag"
$ws2.Range("N25").Value = "Fredag"
$ws2.Range("A26").Value = "Calle"
$ws2.Range("B26").Value = 1
$ws2.Range("G26").Formula = "=SUM(B26:F26)"
$ws2.Range("I26").Value = "Calle"
$ws2.Range("O26").Formula = "=SUM(J26:N26)"
$ws2.Range("A27").Value = "Kim"
$ws2.Range("B27").Value = 1
$ws2.Range("G27").Formula = "=SUM(B27:F27)"
$ws2.Range("I27").Value = "Kim"
$ws2.Range("O27").Formula = "=SUM(J27:N27)"
$ws2.Range("A28").Value = "Nils"
$ws2.Range("B28").Value = 1
$ws2.Range("G28").Formula = "=SUM(B28:F28)"
$ws2.Range("I28").Value = "Nils"
$ws2.Range("O28").Formula = "=SUM(J28:N28)"
$ws2.Range("A29").Value = "Rasmus"
$ws2.Range("B29").Value = 1
$ws2.Range("G29").Formula = "=SUM(B29:F29)"
$ws2.Range("I29").Value = "Rasmus"
$ws2.Range("O29").Formula = "=SUM(J29:N29)"
$ws2.Range("G30").Formula = "=SUM(G26:G29)"
$ws2.Range("O30").Formula = "=SUM(O26:O29)"
$ws2.Range("C33").Value = "Peer-evaluation"
$ws2.Range("K33").Value = "Modular view"
$ws2.Range("B34").Value = "Måndag"
$ws2.Range("C34").Value = "Tisdag"
$ws2.Range("D34").Value = "Onsdag"
$ws2.Range("E34").Value = "Torsdag"
$ws2.Range("F34").Value = "Fredag"
$ws2.Range("J34").Value = "Måndag"
$ws2.Range("K34").Value = "Tisdag"
$ws2.Range("L34").Value = "Onsdag"
$ws2.Range("M34").Value = "Torsdag"
$ws2.Range("N34").Value = "Fredag"
$ws2.Range("A35").Value = "Calle"
$ws2.Range("C35").Value = 3
$ws2.Range("D35").Value = 2
$ws2.Range("G35").Formula = "=SUM(B35:F35)"
$ws2.Range("I35").Value = "Calle"
$ws2.Range("O35").Formula = "=SUM(J35:N35)"
$ws2.Range("A36").Value = "Kim"
$ws2.Range("C36").Value = 3
$ws2.Range("D36").Value = 2
$ws2.Range("G36").Formula = "=SUM(B36:F36)"
$ws2.Range("I36").Value = "Kim"
$ws2.Range("O36").Formula = "=SUM(J36:N36)"
$ws2.Range("A37").Value = "Nils"
$ws2.Range("C37").Value = 3
$ws2.Range("G37").Formula = "=SUM(B37:F37)"
$ws2.Range("I37").Value = "Nils"
$ws2.Range("O37").Formula = "=SUM(J37:N37)"
$ws2.Range("A38").Value = "Rasmus"
$ws2.Range("C38").Value = 3
$ws2.Range("D38").Value = 2
$ws2.Range("G38").Formula = "=SUM(B38:F38)"
$ws2.Range("I38").Value = "Rasmus"
$ws2.Range("O38").Formula = "=SUM(J38:N38)"
$ws2.Range("G39").Formula = "=SUM(G35:G38)"
$ws2.Range("O39").Formula = "=SUM(O35:O38)"
$ws2.Range("C42").Value = "Modular view"
$ws2.Range("K42").Value = "Execution view"
$ws2.Range("B43").Value = "Måndag"
$ws2.Range("C43").Value = "Tisdag"
$ws2.Range("D43").Value = "Onsdag"
$ws2.Range("E43").Value = "Torsdag"
$ws2.Range("F43").Value = "Fredag"
$ws2.Range("J43").Value = "Måndag"
$ws2.Range("K43").Value = "Tisdag"
$ws2.Range("L43").Value = "Onsdag"
$ws2.Range("M43").Value = "Torsdag"
$ws2.Range("N43").Value = "Fredag"
$ws2.Range("A44").Value = "Calle"
$ws2.Range("G44").Formula = "=SUM(B44:F44)"
$ws2.Range("I44").Value = "Calle"
$ws2.Range("O44").Formula = "=SUM(J44:N44)"
$ws2.Range("A45").Value = "Kim"
$ws2.Range("G45").Formula = "=SUM(B45:F45)"
$ws2.Range("I45").Value = "Kim"
$ws2.Range("O45").Formula = "=SUM(J45:N45)"
$ws2.Range("A46").Value = "Nils"
$ws2.Range("G46").Formula = "=SUM(B46:F46)"
$ws2.Range("I46").Value = "Nils"
$ws2.Range("O46").Formula = "=SUM(J46:N46)"
$ws2.Range("A47").Value = "Rasmus"
$ws2.Range("G47").Formula = "=SUM(B47:F47)"
$ws2.Range("I47").Value = "Rasmus"
$ws2.Range("O47").Formula = "=SUM(J47:N47)"
$ws2.Range("G48").Formula = "=SUM(G44:G47)"
$ws2.Range("O48").Formula = "=SUM(O44:O47)"
$ws2.Range("C51").Value = "Execution view"
$ws2.Range("K51").Value = "Architecture evaluation"
$ws2.Range("B52").Value = "Måndag"
$ws2.Range("C52").Value = "Tisdag"
$ws2.Range("D52").Value = "Onsdag"
$ws2.Range("E52").Value = "Torsdag"
$ws2.Range("F52").Value = "Fredag"
$ws2.Range("J52").Value = "Måndag"
$ws2.Range("K52").Value = "Tisdag"
$ws2.Range("L52").Value = "Onsdag"
$ws2.Range("M52").Value = "Torsdag"
$ws2.Range("N52").Value = "Fredag"
$ws2.Range("A53").Value = "Calle"
$ws2.Range("G53").Formula = "=SUM(B53:F53)"
$ws2.Range("I53").Value = "Calle"
$ws2.Range("O53").Formula = "=SUM(J53:N53)"
$ws2.Range("A54").Value = "Kim"
$ws2.Range("G54").Formula = "=SUM(B54:F54)"
$ws2.Range("I54").Value = "Kim"
$ws2.Range("O54").Formula = "=SUM(J54:N54)"
$ws2.Range("A55").Value = "Nils"
$ws2.Range("G55").Formula = "=SUM(B55:F55)"
$ws2.Range("I55").Value = "Nils"
$ws2.Range("O55").Formula = "=SUM(J55:N55)"
$ws2.Range("A56").Value = "Rasmus"
$ws2.Range("G56").Formula = "=SUM(B56:F56)"
$ws2.Range("I56").Value = "Rasmus"
$ws2.Range("O56").Formula = "=SUM(J56:N56)"
$ws2.Range("G57").Formula = "=SUM(G53:G56)"
$ws2.Range("O57").Formula = "=SUM(O53:O56)"

# --- Styles (copied from analogous donor cells on A01) ---
$ws1.Range("G10").Copy()
$ws2.Range("G9").PasteSpecial(-4122)
$ws2.Range("O9").PasteSpecial(-4122)
$ws2.Range("G16").PasteSpecial(-4122)
$ws2.Range("O16").PasteSpecial(-4122)
$ws2.Range("G23").PasteSpecial(-4122)
$ws2.Range("O23").PasteSpecial(-4122)
$ws2.Range("G30").PasteSpecial(-4122)
$ws2.Range("O30").PasteSpecial(-4122)
$ws2.Range("G39").PasteSpecial(-4122)
$ws2.Range("O39").PasteSpecial(-4122)
$ws2.Range("G48").PasteSpecial(-4122)
$ws2.Range("O48").PasteSpecial(-4122)
$ws2.Range("G57").PasteSpecial(-4122)
$ws2.Range("O57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B7").Copy()
$ws2.Range("K19").PasteSpecial(-4122)
$ws2.Range("K22").PasteSpecial(-4122)
$ws2.Range("B26").PasteSpecial(-4122)
$ws2.Range("B27").PasteSpecial(-4122)
$ws2.Range("B28").PasteSpecial(-4122)
$ws2.Range("B29").PasteSpecial(-4122)
$ws2.Range("C35").PasteSpecial(-4122)
$ws2.Range("D35").PasteSpecial(-4122)
$ws2.Range("C36").PasteSpecial(-4122)
$ws2.Range("D36").PasteSpecial(-4122)
$ws2.Range("C37").PasteSpecial(-4122)
$ws2.Range("C38").PasteSpecial(-4122)
$ws2.Range("D38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B36").Copy()
$ws2.Range("B6").PasteSpecial(-4122)
$ws2.Range("J6").PasteSpecial(-4122)
$ws2.Range("J26").PasteSpecial(-4122)
$ws2.Range("J27").PasteSpecial(-4122)
$ws2.Range("J28").PasteSpecial(-4122)
$ws2.Range("J29").PasteSpecial(-4122)
$ws2.Range("B35").PasteSpecial(-4122)
$ws2.Range("J35").PasteSpecial(-4122)
$ws2.Range("K35").PasteSpecial(-4122)
$ws2.Range("B36").PasteSpecial(-4122)
$ws2.Range("J36").PasteSpecial(-4122)
$ws2.Range("K36").PasteSpecial(-4122)
$ws2.Range("B37").PasteSpecial(-4122)
$ws2.Range("J37").PasteSpecial(-4122)
$ws2.Range("K37").PasteSpecial(-4122)
$ws2.Range("B38").PasteSpecial(-4122)
$ws2.Range("J38").PasteSpecial(-4122)
$ws2.Range("K38").PasteSpecial(-4122)
$ws2.Range("B39").PasteSpecial(-4122)
$ws2.Range("J39").PasteSpecial(-4122)
$ws2.Range("B44").PasteSpecial(-4122)
$ws2.Range("C44").PasteSpecial(-4122)
$ws2.Range("J44").PasteSpecial(-4122)
$ws2.Range("K44").PasteSpecial(-4122)
$ws2.Range("B45").PasteSpecial(-4122)
$ws2.Range("C45").PasteSpecial(-4122)
$ws2.Range("J45").PasteSpecial(-4122)
$ws2.Range("K45").PasteSpecial(-4122)
$ws2.Range("B46").PasteSpecial(-4122)
$ws2.Range("C46").PasteSpecial(-4122)
$ws2.Range("J46").PasteSpecial(-4122)
$ws2.Range("K46").PasteSpecial(-4122)
$ws2.Range("B47").PasteSpecial(-4122)
$ws2.Range("C47").PasteSpecial(-4122)
$ws2.Range("J47").PasteSpecial(-4122)
$ws2.Range("K47").PasteSpecial(-4122)
$ws2.Range("B48").PasteSpecial(-4122)
$ws2.Range("J48").PasteSpecial(-4122)
$ws2.Range("B53").PasteSpecial(-4122)
$ws2.Range("C53").PasteSpecial(-4122)
$ws2.Range("J53").PasteSpecial(-4122)
$ws2.Range("K53").PasteSpecial(-4122)
$ws2.Range("B54").PasteSpecial(-4122)
$ws2.Range("C54").PasteSpecial(-4122)
$ws2.Range("J54").PasteSpecial(-4122)
$ws2.Range("K54").PasteSpecial(-4122)
$ws2.Range("B55").PasteSpecial(-4122)
$ws2.Range("C55").PasteSpecial(-4122)
$ws2.Range("J55").PasteSpecial(-4122)
$ws2.Range("K55").PasteSpecial(-4122)
$ws2.Range("B56").PasteSpecial(-4122)
$ws2.Range("C56").PasteSpecial(-4122)
$ws2.Range("J56").PasteSpecial(-4122)
$ws2.Range("K56").PasteSpecial(-4122)
$ws2.Range("B57").PasteSpecial(-4122)
$ws2.Range("J57").PasteSpecial(-4122)
$excel.CutCopyMode = $false


# --- Sheet1 (A01) view state update ---
$ws1.Activate()
$ws1.Range("A2:O40").Select()
